# Auto-generated edit script: updates Leve market-price columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across all 8 job sheets, per the scheduled market-data refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1333.5883
$ws.Range("I15").Value = 1333.5883
$ws.Range("K15").Value = 4000.7649
$ws.Range("M15").Value = -3831.7649
$ws.Range("H33").Value = 275.46667
$ws.Range("J33").Value = 762.6667
$ws.Range("L33").Value = 762.6667
$ws.Range("N33").Value = -1220.6667
$ws.Range("H116").Value = 26433.637
$ws.Range("I116").Value = 4837.3335
$ws.Range("K116").Value = 4837.3335
$ws.Range("M116").Value = -1395.3335
$ws.Range("H137").Value = 17681.273
$ws.Range("I137").Value = 4499.3335
$ws.Range("K137").Value = 13498.0005
$ws.Range("M137").Value = -10948.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 594.5
$ws.Range("I2").Value = 466.72726
$ws.Range("K2").Value = 466.72726
$ws.Range("M2").Value = -353.72726
$ws.Range("H45").Value = 1925.6428
$ws.Range("I45").Value = 1895.9
$ws.Range("K45").Value = 1895.9
$ws.Range("M45").Value = -1518.9
$ws.Range("H61").Value = 12756.833
$ws.Range("I61").Value = 17138
$ws.Range("J61").Value = 3994.5
$ws.Range("K61").Value = 17138
$ws.Range("L61").Value = 3994.5
$ws.Range("M61").Value = -16926
$ws.Range("N61").Value = -4418.5
$ws.Range("H102").Value = 9743.066000000001
$ws.Range("I102").Value = 9743.066000000001
$ws.Range("K102").Value = 9743.066000000001
$ws.Range("M102").Value = -8121.066000000001
$ws.Range("H110").Value = 1565.7142
$ws.Range("I110").Value = 1240
$ws.Range("K110").Value = 1240
$ws.Range("M110").Value = 805
$ws.Range("H116").Value = 594.5
$ws.Range("I116").Value = 466.72726
$ws.Range("K116").Value = 466.72726
$ws.Range("M116").Value = 1827.27274
$ws.Range("H132").Value = 1788353.1
$ws.Range("I132").Value = 2085808.2
$ws.Range("K132").Value = 6257424.6
$ws.Range("M132").Value = -6254894.6
$ws.Range("H136").Value = 12756.833
$ws.Range("I136").Value = 17138
$ws.Range("J136").Value = 3994.5
$ws.Range("K136").Value = 51414
$ws.Range("L136").Value = 11983.5
$ws.Range("M136").Value = -48864
$ws.Range("N136").Value = -17083.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 594.5
$ws.Range("I3").Value = 466.72726
$ws.Range("K3").Value = 466.72726
$ws.Range("M3").Value = -352.72726
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H11").Value = 218.6
$ws.Range("I11").Value = 23.5
$ws.Range("J11").Value = 999
$ws.Range("K11").Value = 23.5
$ws.Range("L11").Value = 999
$ws.Range("M11").Value = 116.5
$ws.Range("N11").Value = -1279
$ws.Range("H12").Value = 289.57144
$ws.Range("J12").Value = 402
$ws.Range("L12").Value = 402
$ws.Range("N12").Value = -738
$ws.Range("H86").Value = 3093.7144
$ws.Range("I86").Value = 3697.5
$ws.Range("J86").Value = 2288.6667
$ws.Range("K86").Value = 3697.5
$ws.Range("L86").Value = 2288.6667
$ws.Range("M86").Value = -2574.5
$ws.Range("N86").Value = -4534.6667
$ws.Range("H89").Value = 3093.7144
$ws.Range("I89").Value = 3697.5
$ws.Range("J89").Value = 2288.6667
$ws.Range("K89").Value = 18487.5
$ws.Range("L89").Value = 11443.3335
$ws.Range("M89").Value = -12871.5
$ws.Range("N89").Value = -22675.3335
$ws.Range("H99").Value = 8311
$ws.Range("J99").Value = 1547.5
$ws.Range("L99").Value = 1547.5
$ws.Range("N99").Value = -4543.5
$ws.Range("H105").Value = 3227
$ws.Range("I105").Value = 2526.9333
$ws.Range("K105").Value = 2526.9333
$ws.Range("M105").Value = -779.9333000000001
$ws.Range("H107").Value = 15626031
$ws.Range("I107").Value = 17858198
$ws.Range("J107").Value = 864.5
$ws.Range("K107").Value = 17858198
$ws.Range("L107").Value = 864.5
$ws.Range("M107").Value = -17856278
$ws.Range("N107").Value = -4704.5
$ws.Range("H134").Value = 9408.040000000001
$ws.Range("I134").Value = 5633.375
$ws.Range("K134").Value = 16900.125
$ws.Range("M134").Value = -14365.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 1475
$ws.Range("I21").Value = 1475
$ws.Range("K21").Value = 1475
$ws.Range("M21").Value = -1240
$ws.Range("H22").Value = 1741.6923
$ws.Range("I22").Value = 275
$ws.Range("J22").Value = 2998.8572
$ws.Range("K22").Value = 275
$ws.Range("L22").Value = 2998.8572
$ws.Range("M22").Value = 75
$ws.Range("N22").Value = -3698.8572
$ws.Range("H25").Value = 3880
$ws.Range("I25").Value = 3880
$ws.Range("K25").Value = 3880
$ws.Range("M25").Value = -3706
$ws.Range("H31").Value = 3735.5278
$ws.Range("J31").Value = 2936
$ws.Range("L31").Value = 2936
$ws.Range("N31").Value = -3526
$ws.Range("H34").Value = 3735.5278
$ws.Range("J34").Value = 2936
$ws.Range("L34").Value = 2936
$ws.Range("N34").Value = -3340
$ws.Range("H58").Value = 15343.182
$ws.Range("J58").Value = 25516.6
$ws.Range("L58").Value = 25516.6
$ws.Range("N58").Value = -25922.6
$ws.Range("H99").Value = 20596.334
$ws.Range("I99").Value = 22105.092
$ws.Range("K99").Value = 22105.092
$ws.Range("M99").Value = -20607.092
$ws.Range("H105").Value = 8076.625
$ws.Range("I105").Value = 10975.9
$ws.Range("K105").Value = 10975.9
$ws.Range("M105").Value = -9228.9
$ws.Range("H126").Value = 20596.334
$ws.Range("I126").Value = 22105.092
$ws.Range("K126").Value = 66315.276
$ws.Range("M126").Value = -63845.276
$ws.Range("H136").Value = 15343.182
$ws.Range("J136").Value = 25516.6
$ws.Range("L136").Value = 76549.79999999999
$ws.Range("N136").Value = -81649.79999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10721.777
$ws.Range("I3").Value = 3832.3333
$ws.Range("J3").Value = 14166.5
$ws.Range("K3").Value = 11496.9999
$ws.Range("L3").Value = 42499.5
$ws.Range("M3").Value = -11384.9999
$ws.Range("N3").Value = -42723.5
$ws.Range("H26").Value = 307.9
$ws.Range("I26").Value = 50
$ws.Range("K26").Value = 150
$ws.Range("M26").Value = 138
$ws.Range("H80").Value = 9684.691999999999
$ws.Range("I80").Value = 4999.3335
$ws.Range("J80").Value = 11090.3
$ws.Range("K80").Value = 14998.0005
$ws.Range("L80").Value = 33270.89999999999
$ws.Range("M80").Value = -14062.0005
$ws.Range("N80").Value = -35142.89999999999
$ws.Range("H83").Value = 9684.691999999999
$ws.Range("I83").Value = 4999.3335
$ws.Range("J83").Value = 11090.3
$ws.Range("K83").Value = 44994.0015
$ws.Range("L83").Value = 99812.7
$ws.Range("M83").Value = -40314.0015
$ws.Range("N83").Value = -109172.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2856.0344
$ws.Range("I113").Value = 2124.0435
$ws.Range("K113").Value = 2124.0435
$ws.Range("M113").Value = 45.95649999999978
$ws.Range("H122").Value = 8657.333000000001
$ws.Range("I122").Value = 2979.75
$ws.Range("K122").Value = 8939.25
$ws.Range("M122").Value = -6489.25
$ws.Range("H132").Value = 7270.558
$ws.Range("I132").Value = 5526.647
$ws.Range("K132").Value = 16579.941
$ws.Range("M132").Value = -14049.941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 924
$ws.Range("I16").Value = 781.4286
$ws.Range("K16").Value = 781.4286
$ws.Range("M16").Value = -611.4286
$ws.Range("H61").Value = 11447.723
$ws.Range("I61").Value = 9827.117
$ws.Range("K61").Value = 9827.117
$ws.Range("M61").Value = -9625.117
$ws.Range("H100").Value = 5500
$ws.Range("I100").Value = 5500
$ws.Range("K100").Value = 5500
$ws.Range("M100").Value = -4959
$ws.Range("H113").Value = 11447.723
$ws.Range("I113").Value = 9827.117
$ws.Range("K113").Value = 9827.117
$ws.Range("M113").Value = -7657.117
$ws.Range("H132").Value = 2392447.8
$ws.Range("I132").Value = 3862402.5
$ws.Range("K132").Value = 11587207.5
$ws.Range("M132").Value = -11584677.5
$ws.Range("H136").Value = 12179.154
$ws.Range("I136").Value = 6165.1665
$ws.Range("K136").Value = 18495.4995
$ws.Range("M136").Value = -15945.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1801.6666
$ws.Range("I107").Value = 1545.0883
$ws.Range("J107").Value = 2424.7856
$ws.Range("K107").Value = 4635.2649
$ws.Range("L107").Value = 7274.3568
$ws.Range("M107").Value = -2715.2649
$ws.Range("N107").Value = -11114.3568
$ws.Range("H122").Value = 282954
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H126").Value = 3663.25
$ws.Range("I126").Value = 3801.4
$ws.Range("K126").Value = 11404.2
$ws.Range("M126").Value = -8934.200000000001
$ws.Range("H136").Value = 1377.7273
$ws.Range("I136").Value = 1100.1428
$ws.Range("K136").Value = 3300.4284
$ws.Range("M136").Value = -750.4284000000002
